# Update the crawl timestamp on every data row (rows 2-404, column O)
# from "2022-12-22 12:56:18" to "2022-12-22 20:49:38", and refresh the
# rating figures that were re-crawled for row 296 (id 5853824).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 404; $row++) {
    $ws.Range("O$row").Value = "2022-12-22 20:49:38"
}

$ws.Range("D296").Value = 4
$ws.Range("E296").Value = 4
